# Generate Report for Handback
# Adds a new handback entry (58fa31bb-1f62-4601-87d4-d0da27664ce6) as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the layout
# already used by the existing rows (188bfde0-... / 8b1927f0-...).

$wb = $excel.ActiveWorkbook

$mdName   = "58fa31bb-1f62-4601-87d4-d0da27664ce6.md"
$zhXlf    = "58fa31bb-1f62-4601-87d4-d0da27664ce6.83615fb1b37cf53a5cff7f0f6c9e441dca1ef6f4.zh-cn.xlf"
$deXlf    = "58fa31bb-1f62-4601-87d4-d0da27664ce6.83615fb1b37cf53a5cff7f0f6c9e441dca1ef6f4.de-de.xlf"
$statusOk = "Handed back: in sync with en-US"

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/58fa31bb-1f62-4601-87d4-d0da27664ce6.md"

# ---------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B4").Value = $statusOk
$wsOverview.Range("C4").Value = $statusOk
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $mdUrl, "", "", $mdName) | Out-Null

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhMdUrl  = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/HEAD/e2e/58fa31bb-1f62-4601-87d4-d0da27664ce6.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/HEAD/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/58fa31bb-1f62-4601-87d4-d0da27664ce6.83615fb1b37cf53a5cff7f0f6c9e441dca1ef6f4.zh-cn.xlf"

$wsZh.Range("B4").Value = $statusOk
$wsZh.Range("D4").Value = "2016-01-27 08:04:42"
$wsZh.Range("G4").Value = "2016-01-27 08:05:26"
$wsZh.Range("H4").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $zhMdUrl, "", "", $mdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), $zhXlfUrl, "", "", $zhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E4"), $zhMdUrl, "", "", $mdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), $zhXlfUrl, "", "", $zhXlf) | Out-Null

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deMdUrl  = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/HEAD/e2e/58fa31bb-1f62-4601-87d4-d0da27664ce6.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/HEAD/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/58fa31bb-1f62-4601-87d4-d0da27664ce6.83615fb1b37cf53a5cff7f0f6c9e441dca1ef6f4.de-de.xlf"

$wsDe.Range("B4").Value = $statusOk
$wsDe.Range("D4").Value = "2016-01-27 08:04:53"
$wsDe.Range("G4").Value = "2016-01-27 08:05:46"
$wsDe.Range("H4").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $deMdUrl, "", "", $mdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), $deXlfUrl, "", "", $deXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E4"), $deMdUrl, "", "", $mdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), $deXlfUrl, "", "", $deXlf) | Out-Null
